$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.630.49"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "2.113.93"
$ws.Range("E3").Value = "  +9.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'254.39"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").Value = "'0.664"
$ws.Range("E6").Value = "  -5.46%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'47.30"
$ws.Range("E8").Value = "  +6.45%  "
$ws.Range("D9").Value = "'60.63"
$ws.Range("E9").Value = "  +3.37%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "'0.0746"
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "2.418.54"
$ws.Range("E13").Value = "  +9.69%  "
$ws.Range("D14").Value = "'14.27"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "'0.832"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "2.108.90"
$ws.Range("E16").Value = "  +9.43%  "
$ws.Range("D17").Value = "'5.12"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "36.621.28"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'73.47"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("D21").Value = "'13.23"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").Value = "'240.49"
$ws.Range("E22").Value = "  -4.66%  "
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'2.48"
$ws.Range("E25").Value = "  -7.48%  "
$ws.Range("D26").Value = "'172.36"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").Value = "'21.62"
$ws.Range("E27").Value = "  +13.97%  "
$ws.Range("E28").Value = "  +3.47%  "
$ws.Range("E29").Value = "  -9.21%  "
$ws.Range("D30").Value = "'29.66"
$ws.Range("E30").Value = "  +65.53%  "
$ws.Range("E31").Value = "  -4.95%  "
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0600"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.0937"
$ws.Range("E34").Value = "  +6.83%  "
$ws.Range("D35").Value = "'0.959"
$ws.Range("E35").Value = "  +6.68%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  +14.48%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.89"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'4.10"
$ws.Range("E39").Value = "  -6.47%  "
$ws.Range("D40").Value = "'1.33"
$ws.Range("E40").Value = "  -11.47%  "
$ws.Range("E41").Value = "  +6.24%  "
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "'98.63"
$ws.Range("E43").Value = "  -8.04%  "
$ws.Range("D44").Value = "'2.78"
$ws.Range("E44").Value = "  +7.77%  "
$ws.Range("D45").Value = "'15.94"
$ws.Range("E45").Value = "  -8.34%  "
$ws.Range("D46").Value = "1.349.30"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0842"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.16"
$ws.Range("E48").Value = "  +11.18%  "
$ws.Range("D49").Value = "2.297.01"
$ws.Range("E49").Value = "  +9.21%  "
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").Value = "'2.29"
$ws.Range("E51").Value = "  -4.70%  "
